$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = '[{''audio_filename'': ''2023-01-15_14-30-25.wav'', ''audio_path'': ''data/audio_files/2023-01-15_14-30-25.wav'', ''audio_datetime'': Timestamp(''2023-01-15 14:30:25''), ''time_diff_minutes'': 4.916666666666667, ''match_score'': 0.8762500000000001, ''extracted_data'': {''age'': 25, ''sex'': ''M'', ''mechanism'': ''motor vehicle accident'', ''injuries'': ''Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture'', ''transcript'': ''This is a 25-year-old male patient involved in a motor vehicle accident. Patient was unconscious at the scene. Head injury with possible skull fracture. Chest trauma with rib fractures. Patient is currently intubated and on mechanical ventilation. Blood pressure is stable. Heart rate is elevated at 110. Patient needs immediate trauma activation.'', ''activation_page'': ''Trauma activation for 25-year-old M with motor vehicle accident. Injuries: Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture. Trauma team activation required.''}}]'
$ws.Range("K2").Value = '{''audio_filename'': ''2023-01-15_14-30-25.wav'', ''audio_path'': ''data/audio_files/2023-01-15_14-30-25.wav'', ''audio_datetime'': Timestamp(''2023-01-15 14:30:25''), ''time_diff_minutes'': 4.916666666666667, ''match_score'': 0.8762500000000001, ''extracted_data'': {''age'': 25, ''sex'': ''M'', ''mechanism'': ''motor vehicle accident'', ''injuries'': ''Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture'', ''transcript'': ''This is a 25-year-old male patient involved in a motor vehicle accident. Patient was unconscious at the scene. Head injury with possible skull fracture. Chest trauma with rib fractures. Patient is currently intubated and on mechanical ventilation. Blood pressure is stable. Heart rate is elevated at 110. Patient needs immediate trauma activation.'', ''activation_page'': ''Trauma activation for 25-year-old M with motor vehicle accident. Injuries: Head Injury, Chest Trauma, Fracture, Rib Fracture, Skull Fracture. Trauma team activation required.''}}'
$ws.Range("J3").Value = '[{''audio_filename'': ''2023-01-15_15-15-30.wav'', ''audio_path'': ''data/audio_files/2023-01-15_15-15-30.wav'', ''audio_datetime'': Timestamp(''2023-01-15 15:15:30''), ''time_diff_minutes'': 4.75, ''match_score'': 0.8620833333333333, ''extracted_data'': {''age'': 32, ''sex'': ''F'', ''mechanism'': ''fall from height'', ''injuries'': ''Fracture, Internal Bleeding'', ''transcript'': ''32-year-old female patient fell from height, approximately 20 feet. Multiple fractures including femur, pelvis, and spine. Internal bleeding suspected. Patient is conscious but in severe pain. Blood pressure is dropping. Heart rate is 130. Need immediate trauma activation for this high fall victim.'', ''activation_page'': ''Trauma activation for 32-year-old F with fall from height. Injuries: Fracture, Internal Bleeding. Trauma team activation required.''}}]'
$ws.Range("K3").Value = '{''audio_filename'': ''2023-01-15_15-15-30.wav'', ''audio_path'': ''data/audio_files/2023-01-15_15-15-30.wav'', ''audio_datetime'': Timestamp(''2023-01-15 15:15:30''), ''time_diff_minutes'': 4.75, ''match_score'': 0.8620833333333333, ''extracted_data'': {''age'': 32, ''sex'': ''F'', ''mechanism'': ''fall from height'', ''injuries'': ''Fracture, Internal Bleeding'', ''transcript'': ''32-year-old female patient fell from height, approximately 20 feet. Multiple fractures including femur, pelvis, and spine. Internal bleeding suspected. Patient is conscious but in severe pain. Blood pressure is dropping. Heart rate is 130. Need immediate trauma activation for this high fall victim.'', ''activation_page'': ''Trauma activation for 32-year-old F with fall from height. Injuries: Fracture, Internal Bleeding. Trauma team activation required.''}}'
